$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append below the existing data (row 21 is the last used row).
$row = 22

# --- Values -----------------------------------------------------------
$ws.Cells.Item($row, 1).Value = 61186
$ws.Cells.Item($row, 2).Value = 4
$ws.Cells.Item($row, 3).Value = "Nutriflex Omega special, Infusionsemulsion 2500 ml"
$ws.Cells.Item($row, 4).Value = "B. Braun Medical AG"
$ws.Cells.Item($row, 5).Value = "07.01.2."
$ws.Cells.Item($row, 6).Value = "B05BA10"
$ws.Cells.Item($row, 7).Value = "Synthetika human"
$ws.Cells.Item($row, 8).Value = 40522
$ws.Cells.Item($row, 9).Value = 40522
$ws.Cells.Item($row, 10).Value = 42347
$ws.Cells.Item($row, 11).Value = 4
$ws.Cells.Item($row, 12).Value = "5 x 2500 ml"
$ws.Cells.Item($row, 13).Value = "Beutel"
$ws.Cells.Item($row, 14).Value = "B"
$ws.Cells.Item($row, 15).Value = "glucosum anhydricum, natrii dihydrogenophosphas dihydricus, zinci acetas dihydricus, isoleucinum, leucinum, lysinum anhydricum, methioninum, phenylalaninum, threoninum, tryptophanum, valinum, argininum, histidinum, alaninum, acidum asparticum, acidum glutamicum, glycinum, prolinum, serinum, natrii hydroxidum, natrii chloridum, natrii acetas trihydricus, kalii acetas, magnesii acetas tetrahydricus, calcii chloridum dihydricum, aminoacida, nitrogenia, carbohydrata, materia crassa, natrium, kalium, magnesium, calcium, zincum, chloridum, phosphas, acetas, sojae oleum, triglycerida saturata media, omega-3 acidorum triglycerida"
$ws.Cells.Item($row, 16).Value = "I) Glucoselösung: glucosum anhydricum 360 g ut glucosum monohydricum, natrii dihydrogenophosphas dihydricus 6.24 g, zinci acetas dihydricus 17.56 mg, aqua ad iniectabilia q.s. ad solutionem pro 1000 ml.`nII) Fettemulsion: sojae oleum 40 g, triglycerida saturata media 50 g, omega-3 acidorum triglycerida 10 g, glycerolum, lecithinum ex ovo, natrii oleas, antiox.: E 307 100 mg, aqua ad iniectabilia q.s. ad emulsionem pro 500 ml.`nIII) Aminosäurenlösung: isoleucinum 8.21 g, leucinum 10.96 g, lysinum anhydricum 7.95 g ut lysinum monohydricum, methioninum 6.84 g, phenylalaninum 12.29 g, threoninum 6.35 g, tryptophanum 2 g, valinum 9.01 g, argininum 9.45 g, histidinum 4.38 g ut histidini hydrochloridum monohydricum, alaninum 16.98 g, acidum asparticum 5.25 g, acidum glutamicum 12.27 g, glycinum 5.78 g, prolinum 11.9 g, serinum 10.5 g, natrii hydroxidum 2.928 g, natrii chloridum 946 mg, natrii acetas trihydricus 626 mg, kalii acetas 9.222 g, magnesii acetas tetrahydricus 2.274 g, calcii chloridum dihydricum 1.558 g, aqua ad iniectabilia q.s. ad solutionem pro 1000 ml.`n.`nI) et II) et III) corresp.: aminoacida 57.44 g/l, nitrogenia 8 g/l, carbohydrata 144 g/l, materia crassa 40 g/l, natrium 53.6 mmol/l, kalium 37.6 mmol/l, magnesium 4.24 mmol/l, calcium 4.24 mmol/l, zincum 0.032 mmol/l, chloridum 48 mmol/l, phosphas 16 mmol/l, acetas 48 mmol/l, in emulsione recenter mixta 1000 ml.`nCorresp. 4941 kJ pro 1 l."
$ws.Cells.Item($row, 17).Value = "Parenterale Ernährung"
# Column R (18) is left empty (no content), matching the source row, but we
# still touch it so an (empty) cell record exists at R22.
$ws.Cells.Item($row, 18).NumberFormat = "GENERAL"

# --- Formatting ---------------------------------------------------------
# Clone number format / alignment from the row above (row 21, the previous
# last data row) for columns A:Q so every column keeps the sheet's existing
# per-column look (right-most column R stays on the plain column default).
$ws.Range("A21:Q21").Copy()
$ws.Range("A22:Q22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The new "u" / composition column (P) additionally gets wrap-text turned on
# (a new style vs. the rest of the "general text" columns).
$ws.Cells.Item($row, 16).WrapText = $true

$ws.Rows.Item($row).RowHeight = 12.75

$ws.Range("A22").Select()
